# The presentation carries two SharePoint-related custom XML parts (the
# document-library "content type" plumbing): one holds the SharePoint
# "FormTemplates" content-type XML, the other holds the "documentManagement"
# properties XML (Status / MediaServiceKeyPoints). This edit swaps which
# part holds which payload - i.e. the part that used to contain the
# FormTemplates XML now contains the documentManagement XML, and vice
# versa.
#
# Reproduce it through the supported automation surface for a
# presentation's custom XML parts: Presentation.CustomXMLParts (Add /
# Item / Delete / XML), by removing the two existing parts and re-adding
# them with their XML payloads swapped.

$p = $ppt.ActivePresentation

$formTemplatesXml = '<?mso-contentType?><FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"><Display>DocumentLibraryForm</Display><Edit>DocumentLibraryForm</Edit><New>DocumentLibraryForm</New></FormTemplates>'
$documentManagementXml = '<?xml version="1.0" encoding="utf-8"?><p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"><documentManagement><Status xmlns="71af3243-3dd4-4a8d-8c0d-dd76da1f02a5">Not started</Status><MediaServiceKeyPoints xmlns="71af3243-3dd4-4a8d-8c0d-dd76da1f02a5" xsi:nil="true"/></documentManagement></p:properties>'

$parts = $p.CustomXMLParts

# Drop whichever of the two known parts are currently present (by content),
# regardless of their current Item() order.
for ($i = $parts.Count; $i -ge 1; $i--) {
    $existingXml = [string]$parts.Item($i).XML
    if ($existingXml.Contains("FormTemplates") -or $existingXml.Contains("documentManagement")) {
        $parts.Item($i).Delete()
    }
}

# Re-add them with the payloads swapped.
$parts.Add($documentManagementXml) | Out-Null
$parts.Add($formTemplatesXml) | Out-Null
